# Update the K column (column G) values on the active sheet to reflect
# the regenerated save_data (K instead of Strike#, recalculated std/mean).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1
    3  = 1
    4  = 2
    5  = 2
    6  = 0
    7  = 1
    8  = 1
    9  = 1
    10 = 4
    11 = 0
    12 = 0
    13 = 0
    14 = 1
    15 = 1
    16 = 1
    17 = 0
    19 = 2
    20 = 1
    21 = 1
    22 = 1
    23 = 1
    24 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
